$wb = $excel.ActiveWorkbook
$stevens = $wb.Worksheets.Item("Stevens_2012_placebo")

# --- 1. Create "WrongSheet": a copy of Stevens_2012_placebo containing only
#        the header row (used to test that reading fails / is rejected). ---
$stevens.Copy($null, $stevens)
$wrongSheet = $wb.Worksheets.Item($stevens.Index + 1)
$wrongSheet.Name = "WrongSheet"
$wrongSheet.Rows("2:78").Delete()
$wrongSheet.Range("A1:N1").Select()

# --- 2. Create "CorrectSheet_additionalCols": a copy of Stevens_2012_placebo
#        with extra (ignorable) columns appended, to prove that reading
#        parameters tolerates additional columns in the sheet. ---
$stevens.Copy($null, $wrongSheet)
$correctSheet = $wb.Worksheets.Item($wrongSheet.Index + 1)
$correctSheet.Name = "CorrectSheet_additionalCols"
$correctSheet.Range("Q6").Value = 12
$correctSheet.Range("Q7").Value = 2
$correctSheet.Range("Q8").Value = 3
$correctSheet.Range("A1:N78").Select()

# --- 3. Update the selection on Stevens_2012_placebo (no longer the active
#        tab) and make CorrectSheet_additionalCols the active sheet. ---
$stevens.Activate()
$stevens.Range("A1:N78").Select()

$correctSheet.Activate()

$sheetNames = @()
foreach ($s in $wb.Worksheets) { $sheetNames += $s.Name }
Write-Output ($sheetNames -join ", ")
